$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    # Force the cell to stay a text cell (mirrors the inlineStr cells in the
    # source file) even though many of the values look like plain numbers.
    $Cell.Style = "Normal"
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

# --- Row 42 / 43: ranking reshuffle swaps Arweave and Cosmos ---
Set-TextValue $ws.Cells.Item(42, 2) "Cosmos"
Set-TextValue $ws.Cells.Item(42, 3) "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Cells.Item(42, 4) "8.84"
Set-TextValue $ws.Cells.Item(42, 5) "  +5.49%  "

Set-TextValue $ws.Cells.Item(43, 2) "Arweave"
Set-TextValue $ws.Cells.Item(43, 3) "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue $ws.Cells.Item(43, 4) "44.23"
Set-TextValue $ws.Cells.Item(43, 5) "  -7.91%  "

# --- Price (D) / Volume(1h) (E) refresh for the remaining rows ---
Set-TextValue $ws.Cells.Item(2, 4) "69.892.59"
Set-TextValue $ws.Cells.Item(2, 5) "  +1.76%  "
Set-TextValue $ws.Cells.Item(3, 4) "3.730.62"
Set-TextValue $ws.Cells.Item(3, 5) "  +18.61%  "
Set-TextValue $ws.Cells.Item(4, 5) "  +0.00%  "
Set-TextValue $ws.Cells.Item(5, 4) "614.54"
Set-TextValue $ws.Cells.Item(5, 5) "  +6.24%  "
Set-TextValue $ws.Cells.Item(6, 4) "177.32"
Set-TextValue $ws.Cells.Item(6, 5) "  -1.43%  "
Set-TextValue $ws.Cells.Item(7, 4) "3.727.55"
Set-TextValue $ws.Cells.Item(7, 5) "  +18.52%  "
Set-TextValue $ws.Cells.Item(8, 5) "  +0.14%  "
Set-TextValue $ws.Cells.Item(9, 4) "0.541"
Set-TextValue $ws.Cells.Item(9, 5) "  +3.42%  "
Set-TextValue $ws.Cells.Item(10, 4) "0.167"
Set-TextValue $ws.Cells.Item(10, 5) "  +9.46%  "
Set-TextValue $ws.Cells.Item(11, 4) "6.39"
Set-TextValue $ws.Cells.Item(11, 5) "  -1.91%  "
Set-TextValue $ws.Cells.Item(12, 4) "0.500"
Set-TextValue $ws.Cells.Item(12, 5) "  +6.29%  "
Set-TextValue $ws.Cells.Item(13, 4) "40.77"
Set-TextValue $ws.Cells.Item(13, 5) "  +10.29%  "
Set-TextValue $ws.Cells.Item(14, 5) "  +5.47%  "
Set-TextValue $ws.Cells.Item(15, 4) "4.352.36"
Set-TextValue $ws.Cells.Item(15, 5) "  +18.62%  "
Set-TextValue $ws.Cells.Item(16, 4) "3.732.54"
Set-TextValue $ws.Cells.Item(16, 5) "  +18.68%  "
Set-TextValue $ws.Cells.Item(17, 4) "69.901.83"
Set-TextValue $ws.Cells.Item(17, 5) "  +1.93%  "
Set-TextValue $ws.Cells.Item(18, 5) "  +0.96%  "
Set-TextValue $ws.Cells.Item(19, 5) "  +5.99%  "
Set-TextValue $ws.Cells.Item(20, 4) "515.74"
Set-TextValue $ws.Cells.Item(20, 5) "  +5.40%  "
Set-TextValue $ws.Cells.Item(21, 4) "16.70"
Set-TextValue $ws.Cells.Item(21, 5) "  +1.40%  "
Set-TextValue $ws.Cells.Item(22, 4) "9.37"
Set-TextValue $ws.Cells.Item(22, 5) "  +20.23%  "
Set-TextValue $ws.Cells.Item(23, 4) "0.728"
Set-TextValue $ws.Cells.Item(23, 5) "  +4.20%  "
Set-TextValue $ws.Cells.Item(24, 4) "88.33"
Set-TextValue $ws.Cells.Item(24, 5) "  +5.08%  "
Set-TextValue $ws.Cells.Item(25, 5) "  +5.88%  "
Set-TextValue $ws.Cells.Item(26, 4) "13.57"
Set-TextValue $ws.Cells.Item(26, 5) "  +4.22%  "
Set-TextValue $ws.Cells.Item(27, 4) "10.94"
Set-TextValue $ws.Cells.Item(27, 5) "  +3.18%  "
Set-TextValue $ws.Cells.Item(28, 5) "  -0.08%  "
Set-TextValue $ws.Cells.Item(29, 5) "  +33.63%  "
Set-TextValue $ws.Cells.Item(30, 4) "2.51"
Set-TextValue $ws.Cells.Item(30, 5) "  +6.24%  "
Set-TextValue $ws.Cells.Item(31, 5) "  +7.79%  "
Set-TextValue $ws.Cells.Item(32, 5) "  -3.45%  "
Set-TextValue $ws.Cells.Item(33, 4) "31.39"
Set-TextValue $ws.Cells.Item(33, 5) "  +11.26%  "
Set-TextValue $ws.Cells.Item(34, 5) "  +2.73%  "
Set-TextValue $ws.Cells.Item(35, 4) "1.00"
Set-TextValue $ws.Cells.Item(35, 5) "  -0.04%  "
Set-TextValue $ws.Cells.Item(36, 5) "  +7.58%  "
Set-TextValue $ws.Cells.Item(37, 4) "1.03"
Set-TextValue $ws.Cells.Item(37, 5) "  +7.92%  "
Set-TextValue $ws.Cells.Item(38, 5) "  +5.23%  "
Set-TextValue $ws.Cells.Item(39, 4) "2.18"
Set-TextValue $ws.Cells.Item(39, 5) "  +6.72%  "
Set-TextValue $ws.Cells.Item(40, 4) "0.133"
Set-TextValue $ws.Cells.Item(40, 5) "  +5.35%  "
Set-TextValue $ws.Cells.Item(41, 4) "51.26"
Set-TextValue $ws.Cells.Item(41, 5) "  +4.08%  "
Set-TextValue $ws.Cells.Item(44, 4) "421.37"
Set-TextValue $ws.Cells.Item(44, 5) "  +4.88%  "
Set-TextValue $ws.Cells.Item(45, 4) "3.073.16"
Set-TextValue $ws.Cells.Item(45, 5) "  +9.21%  "
Set-TextValue $ws.Cells.Item(46, 5) "  +0.18%  "
Set-TextValue $ws.Cells.Item(47, 5) "  +4.81%  "
Set-TextValue $ws.Cells.Item(48, 4) "27.89"
Set-TextValue $ws.Cells.Item(48, 5) "  +0.78%  "
Set-TextValue $ws.Cells.Item(49, 4) "2.52"
Set-TextValue $ws.Cells.Item(49, 5) "  +5.79%  "
Set-TextValue $ws.Cells.Item(50, 4) "135.86"
Set-TextValue $ws.Cells.Item(50, 5) "  +0.31%  "
